$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing used range (A1:B7) so we can rebuild the table in its new layout.
$ws.Range("A1:B16").Clear()

# Header row
$ws.Range("A1").Value = "flag"
$ws.Range("B1").Value = "description"

# Rows applicable to all sensors
$ws.Range("A2").Value = "R"
$ws.Range("B2").Value = "removed buoy"

$ws.Range("A3").Value = "D"
$ws.Range("B3").Value = "buoy deployed for season"

# flag_temp section
$ws.Range("A5").Value = "flag_temp"

$ws.Range("A6").Value = "e"
$ws.Range("B6").Value = "data errant, recoded to na"

$ws.Range("A7").Value = "i"
$ws.Range("B7").Value = "data reporting intermittently"

# flag_do section
$ws.Range("A9").Value = "flag_do"

$ws.Range("A10").Value = "w"
$ws.Range("B10").Value = "do sensor cleaned"

$ws.Range("A11").Value = "c"
$ws.Range("B11").Value = "do sensor calibrated"

# flag_met section
$ws.Range("A14").Value = "flag_met"

$ws.Range("A15").Value = "n"
$ws.Range("B15").Value = "weather station adjusted to orient more northward"

$ws.Range("A16").Value = "m"
$ws.Range("B16").Value = "rain gague malfunction; values incorrect - recoded to na"

# Update selection to match the authored state
$ws.Range("C13").Select()
